# Read data from multiple sheets and merge:
# Duplicate the "Staffs" sheet into a new "Staffs_2015" sheet (placed right
# after "Staffs" and before "SalaryLevel"), update it with a new set of
# staff rows (Bill / Ben) and their hyperlinked e-mail addresses, then tidy
# up the selection shown on the "Staffs" sheet.

$wb = $excel.ActiveWorkbook
$staffs = $wb.Worksheets.Item(1)

# --- 1. Create the new worksheet right after "Staffs" -----------------
$staffs2015 = $wb.Worksheets.Add($null, $staffs)
$staffs2015.Name = "Staffs_2015"

# --- 2. Bring over the header / label rows (values + formatting) ------
$staffs.Range("A1:G4").Copy()
$staffs2015.Range("A1").PasteSpecial()
$staffs2015.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

# --- 3. New staff data for 2015 ----------------------------------------
$staffs2015.Range("A5").Value = "Bill"
$staffs2015.Range("C5").Value = 35
$staffs2015.Range("D5").Value = 2
$staffs2015.Range("E5").Value = 3
$staffs2015.Range("F5").Value = 123
$staffs2015.Range("G5").Value = "8A9B1"

$staffs2015.Range("A6").Value = "Ben"
$staffs2015.Range("C6").Value = 65
$staffs2015.Range("D6").Value = 5
$staffs2015.Range("E6").Value = 4
$staffs2015.Range("F6").Value = 234
$staffs2015.Range("G6").Value = "10AB83"

# --- 4. Emails as hyperlinks, styled the same as on "Staffs" -----------
$staffs2015.Range("B5").Value = "bill@gmail.com"
$staffs2015.Hyperlinks.Add($staffs2015.Range("B5"), "mailto:bill@gmail.com")
$staffs2015.Range("B5").Style = "Hyperlink"

$staffs2015.Range("B6").Value = "ben@yahoo.com"
$staffs2015.Hyperlinks.Add($staffs2015.Range("B6"), "mailto:ben@yahoo.com")
$staffs2015.Range("B6").Style = "Hyperlink"

# --- 5. Update selection on "Staffs" (no longer the active sheet) -----
$staffs.Range("A1:G6").Select()

# --- 6. Current selection on the new sheet + make it the active tab ---
$staffs2015.Range("D5").Select()

Write-Output "Staffs_2015 sheet created and populated"
